# Scheduled-runner style refresh of market/profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2681
$ws.Range("I62").Value = 2834.6667
$ws.Range("J62").Value = 2220
$ws.Range("K62").Value = 2834.6667
$ws.Range("L62").Value = 2220
$ws.Range("M62").Value = -2210.6667
$ws.Range("N62").Value = -3468

$ws.Range("H65").Value = 2681
$ws.Range("I65").Value = 2834.6667
$ws.Range("J65").Value = 2220
$ws.Range("K65").Value = 14173.3335
$ws.Range("L65").Value = 11100
$ws.Range("M65").Value = -11053.3335
$ws.Range("N65").Value = -17340

$ws.Range("H138").Value = 4056.6843
$ws.Range("I138").Value = 2017.0952
$ws.Range("J138").Value = 6576.1763
$ws.Range("K138").Value = 6051.2856
$ws.Range("L138").Value = 19728.5289
$ws.Range("M138").Value = -911.2856000000002
$ws.Range("N138").Value = -30008.5289

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 260
$ws.Range("I4").Value = 260
$ws.Range("K4").Value = 260
$ws.Range("M4").Value = -144

$ws.Range("H32").Value = 20205.984
$ws.Range("I32").Value = 22086.666
$ws.Range("J32").Value = 10973.546
$ws.Range("K32").Value = 22086.666
$ws.Range("L32").Value = 10973.546
$ws.Range("M32").Value = -21799.666
$ws.Range("N32").Value = -11547.546

$ws.Range("H37").Value = 5900
$ws.Range("I37").Value = 5900
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 5900
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -5627
$ws.Range("N37").ClearContents()

$ws.Range("H44").Value = 333369660
$ws.Range("I44").Value = 39000
$ws.Range("J44").Value = 500035000
$ws.Range("K44").Value = 39000
$ws.Range("L44").Value = 500035000
$ws.Range("M44").Value = -38512
$ws.Range("N44").Value = -500035976

$ws.Range("H74").Value = 5256.7417
$ws.Range("I74").Value = 2102.5862
$ws.Range("K74").Value = 2102.5862
$ws.Range("M74").Value = -1228.5862

$ws.Range("H77").Value = 5256.7417
$ws.Range("I77").Value = 2102.5862
$ws.Range("K77").Value = 10512.931
$ws.Range("M77").Value = -6144.931

$ws.Range("H118").Value = 30479.8
$ws.Range("J118").Value = 30479.8
$ws.Range("L118").Value = 30479.8
$ws.Range("N118").Value = -33793.8

$ws.Range("H132").Value = 4881.154
$ws.Range("I132").Value = 1805.1765
$ws.Range("J132").Value = 7258.0454
$ws.Range("K132").Value = 5415.529500000001
$ws.Range("L132").Value = 21774.1362
$ws.Range("M132").Value = -2885.529500000001
$ws.Range("N132").Value = -26834.1362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 61270
$ws.Range("J55").Value = 61270
$ws.Range("L55").Value = 61270
$ws.Range("N55").Value = -61816

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2072.7795
$ws.Range("I31").Value = 1515.1777
$ws.Range("K31").Value = 1515.1777
$ws.Range("M31").Value = -1220.1777

$ws.Range("H34").Value = 2072.7795
$ws.Range("I34").Value = 1515.1777
$ws.Range("K34").Value = 1515.1777
$ws.Range("M34").Value = -1313.1777

$ws.Range("H99").Value = 4900
$ws.Range("I99").Value = 4700
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 4700
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -3202
$ws.Range("N99").Value = -7996

$ws.Range("H126").Value = 4900
$ws.Range("I126").Value = 4700
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 14100
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -11630
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 7576612
$ws.Range("J5").Value = 20834872
$ws.Range("L5").Value = 62504616
$ws.Range("N5").Value = -62504840

$ws.Range("H122").Value = 762.96
$ws.Range("I122").Value = 419.45456
$ws.Range("J122").Value = 1032.8572
$ws.Range("K122").Value = 3775.09104
$ws.Range("L122").Value = 9295.7148
$ws.Range("M122").Value = -1325.09104
$ws.Range("N122").Value = -14195.7148

$ws.Range("H135").Value = 7576612
$ws.Range("J135").Value = 20834872
$ws.Range("L135").Value = 187513848
$ws.Range("N135").Value = -187518918

$ws.Range("H138").Value = 12650.546
$ws.Range("I138").Value = 26097.5
$ws.Range("J138").Value = 4966.5713
$ws.Range("K138").Value = 78292.5
$ws.Range("L138").Value = 14899.7139
$ws.Range("M138").Value = -73152.5
$ws.Range("N138").Value = -25179.7139

$ws.Range("H139").Value = 1469644.5
$ws.Range("I139").Value = 2819415.2
$ws.Range("J139").Value = 2502.4348
$ws.Range("K139").Value = 8458245.600000001
$ws.Range("L139").Value = 7507.3044
$ws.Range("M139").Value = -8453105.600000001
$ws.Range("N139").Value = -17787.3044

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6288.4614
$ws.Range("I122").Value = 11598.2
$ws.Range("J122").Value = 2969.875
$ws.Range("K122").Value = 34794.60000000001
$ws.Range("L122").Value = 8909.625
$ws.Range("M122").Value = -32344.60000000001
$ws.Range("N122").Value = -13809.625

$ws.Range("H132").Value = 5839.064
$ws.Range("I132").Value = 3916.7908
$ws.Range("J132").Value = 26503.5
$ws.Range("K132").Value = 11750.3724
$ws.Range("L132").Value = 79510.5
$ws.Range("M132").Value = -9220.3724
$ws.Range("N132").Value = -84570.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7134.654
$ws.Range("I7").Value = 5280.067
$ws.Range("J7").Value = 9663.637
$ws.Range("K7").Value = 5280.067
$ws.Range("L7").Value = 9663.637
$ws.Range("M7").Value = -5168.067
$ws.Range("N7").Value = -9887.637

$ws.Range("H55").Value = 275.29166
$ws.Range("I55").Value = 243.64285
$ws.Range("J55").Value = 319.6
$ws.Range("K55").Value = 243.64285
$ws.Range("L55").Value = 319.6
$ws.Range("M55").Value = -70.64285000000001
$ws.Range("N55").Value = -665.6

$ws.Range("H126").Value = 7134.654
$ws.Range("I126").Value = 5280.067
$ws.Range("J126").Value = 9663.637
$ws.Range("K126").Value = 15840.201
$ws.Range("L126").Value = 28990.911
$ws.Range("M126").Value = -13370.201
$ws.Range("N126").Value = -33930.911

$ws.Range("H132").Value = 3298.1619
$ws.Range("I132").Value = 3339.547
$ws.Range("K132").Value = 10018.641
$ws.Range("M132").Value = -7488.641

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 532.375
$ws.Range("I100").Value = 534.75
$ws.Range("J100").Value = 530
$ws.Range("K100").Value = 1069.5
$ws.Range("L100").Value = 1060
$ws.Range("M100").Value = -528.5
$ws.Range("N100").Value = -2142

$ws.Range("H126").Value = 2088.3333
$ws.Range("I126").Value = 2071.4285
$ws.Range("J126").Value = 2147.5
$ws.Range("K126").Value = 6214.2855
$ws.Range("L126").Value = 6442.5
$ws.Range("M126").Value = -3744.2855
$ws.Range("N126").Value = -11382.5

$ws.Range("H132").Value = 1268.2363
$ws.Range("I132").Value = 642.7353
$ws.Range("K132").Value = 1928.2059
$ws.Range("M132").Value = 601.7940999999998
